$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.679.59"
$ws.Range("E2").Value = "  +3.68%  "

$ws.Range("D3").Value = "3.277.95"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.19"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "629.87"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.390"
$ws.Range("E7").Value = "  +21.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.702"
$ws.Range("E8").Value = "  +18.05%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "3.271.42"
$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("E11").Value = "  -3.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.186"
$ws.Range("E12").Value = "  +11.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -4.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.26"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("D15").Value = "3.877.01"
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").Value = "88.336.59"
$ws.Range("E17").Value = "  +3.67%  "

$ws.Range("D18").Value = "3.279.06"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.13"
$ws.Range("E20").Value = "  -3.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.39"

$ws.Range("E22").Value = "  -2.97%  "

$ws.Range("E23").Value = "  +2.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.37"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.31"
$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.20"
$ws.Range("E26").Value = "  -5.26%  "

$ws.Range("D27").Value = "3.434.88"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "77.11"
$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.184"
$ws.Range("E31").Value = "  +11.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "574.47"
$ws.Range("E33").Value = "  -5.69%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "8.89"
$ws.Range("E34").Value = "  -3.90%  "

$ws.Range("E35").Value = "  -9.81%  "

$ws.Range("E36").Value = "  -4.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.08"
$ws.Range("E37").Value = "  +9.60%  "

$ws.Range("E38").Value = "  -8.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.67"
$ws.Range("E39").Value = "  -2.59%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.79"
$ws.Range("E41").Value = "  +2.62%  "

$ws.Range("E42").Value = "  -4.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.02"
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("E44").Value = "  -4.29%  "

$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.34"
$ws.Range("E46").Value = "  -3.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "180.33"
$ws.Range("E47").Value = "  -5.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.88"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("E49").Value = "  -5.29%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.126"
$ws.Range("E50").Value = "  +12.93%  "

$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0673"
$ws.Range("E51").Value = "  +20.03%  "
